$d = $word.ActiveDocument

# 1) Bump the Unity version referenced in the first paragraph:
#    "2019.4.7f1" -> "2019.4.8f1"
$found = $d.Content.Find.Execute(
    "2019.4.7f1", $true, $false, $false, $false, $false,
    $true, 1, $false, "2019.4.8f1", 2)

# 2) Close up the stray space in "mario odissey" -> "marioodissey"
$found = $d.Content.Find.Execute(
    "mario odissey", $true, $false, $false, $false, $false,
    $true, 1, $false, "marioodissey", 2)
